$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Prevent Excel from auto-converting date-like text (Startdatum/Slutdatum columns Y, AA)
# into date serial numbers when we assign string values.
$ws.Range("Y2:Y20").NumberFormat = "@"
$ws.Range("AA2:AA20").NumberFormat = "@"
$ws.Range("Z2:Z20").NumberFormat = "@"
$ws.Range("AB2:AB20").NumberFormat = "@"

$arr = New-Object "object[,]" 19,51

$arr[0,0] = 71934009
$arr[0,1] = 98520
$arr[0,2] = "Ovaliderad"
$arr[0,3] = "LC"
$arr[0,4] = 222498
$arr[0,5] = "Blåsippa"
$arr[0,6] = "Hepatica nobilis"
$arr[0,7] = "Schreb."
$arr[0,8] = ""
$arr[0,9] = ""
$arr[0,10] = ""
$arr[0,11] = ""
$arr[0,12] = ""
$arr[0,13] = ""
$arr[0,14] = ""
$arr[0,15] = "Skogen NÖ om Bålbroskolan, Upl"
$arr[0,16] = 688625.398946143
$arr[0,17] = 6628251.815486348
$arr[0,18] = 10
$arr[0,19] = "Stockholm"
$arr[0,20] = "Norrtälje"
$arr[0,21] = "Uppland"
$arr[0,22] = "Rimbo"
$arr[0,23] = ""
$arr[0,24] = "2018-06-18"
$arr[0,25] = "00:00"
$arr[0,26] = "2018-06-18"
$arr[0,27] = "00:00"
$arr[0,28] = ""
$arr[0,29] = $false
$arr[0,30] = $false
$arr[0,31] = ""
$arr[0,32] = $false
$arr[0,33] = "Skogsmark"
$arr[0,34] = "Blandskog"
$arr[0,35] = ""
$arr[0,36] = ""
$arr[0,37] = ""
$arr[0,38] = ""
$arr[0,39] = ""
$arr[0,40] = ""
$arr[0,41] = ""
$arr[0,42] = ""
$arr[0,43] = ""
$arr[0,44] = ""
$arr[0,45] = ""
$arr[0,46] = ""
$arr[0,47] = ""
$arr[0,48] = "Beata Larsson"
$arr[0,49] = "Beata Larsson"
$arr[0,50] = ""

$arr[1,0] = 71934047
$arr[1,1] = 98520
$arr[1,2] = "Ovaliderad"
$arr[1,3] = "LC"
$arr[1,4] = 222498
$arr[1,5] = "Blåsippa"
$arr[1,6] = "Hepatica nobilis"
$arr[1,7] = "Schreb."
$arr[1,8] = ""
$arr[1,9] = ""
$arr[1,10] = ""
$arr[1,11] = ""
$arr[1,12] = ""
$arr[1,13] = ""
$arr[1,14] = ""
$arr[1,15] = "Skogen NÖ om Bålbroskolan, Upl"
$arr[1,16] = 688522.0690827493
$arr[1,17] = 6628424.789042529
$arr[1,18] = 5
$arr[1,19] = "Stockholm"
$arr[1,20] = "Norrtälje"
$arr[1,21] = "Uppland"
$arr[1,22] = "Rimbo"
$arr[1,23] = ""
$arr[1,24] = "2018-06-18"
$arr[1,25] = "00:00"
$arr[1,26] = "2018-06-18"
$arr[1,27] = "00:00"
$arr[1,28] = ""
$arr[1,29] = $false
$arr[1,30] = $false
$arr[1,31] = ""
$arr[1,32] = $false
$arr[1,33] = "Skogsmark"
$arr[1,34] = "Blandskog"
$arr[1,35] = ""
$arr[1,36] = ""
$arr[1,37] = ""
$arr[1,38] = ""
$arr[1,39] = ""
$arr[1,40] = ""
$arr[1,41] = ""
$arr[1,42] = ""
$arr[1,43] = ""
$arr[1,44] = ""
$arr[1,45] = ""
$arr[1,46] = ""
$arr[1,47] = ""
$arr[1,48] = "Beata Larsson"
$arr[1,49] = "Beata Larsson"
$arr[1,50] = ""

$arr[2,0] = 71955318
$arr[2,1] = 96355
$arr[2,2] = "Ovaliderad"
$arr[2,3] = "LC"
$arr[2,4] = 219862
$arr[2,5] = "Nästrot"
$arr[2,6] = "Neottia nidus-avis"
$arr[2,7] = "(L.) Rich."
$arr[2,8] = ""
$arr[2,9] = ""
$arr[2,10] = ""
$arr[2,11] = ""
$arr[2,12] = ""
$arr[2,13] = ""
$arr[2,14] = ""
$arr[2,15] = "Björkenäs, Gläntan, Upl"
$arr[2,16] = 688361.9018547204
$arr[2,17] = 6628601.946578045
$arr[2,18] = 5
$arr[2,19] = "Stockholm"
$arr[2,20] = "Norrtälje"
$arr[2,21] = "Uppland"
$arr[2,22] = "Rimbo"
$arr[2,23] = ""
$arr[2,24] = "2018-06-20"
$arr[2,25] = "00:00"
$arr[2,26] = "2018-06-20"
$arr[2,27] = "00:00"
$arr[2,28] = ""
$arr[2,29] = $false
$arr[2,30] = $false
$arr[2,31] = ""
$arr[2,32] = $false
$arr[2,33] = "Skogsmark"
$arr[2,34] = "Blandskog"
$arr[2,35] = ""
$arr[2,36] = ""
$arr[2,37] = ""
$arr[2,38] = ""
$arr[2,39] = ""
$arr[2,40] = ""
$arr[2,41] = ""
$arr[2,42] = ""
$arr[2,43] = ""
$arr[2,44] = ""
$arr[2,45] = ""
$arr[2,46] = ""
$arr[2,47] = ""
$arr[2,48] = "Beata Larsson"
$arr[2,49] = "Beata Larsson"
$arr[2,50] = ""

$arr[3,0] = 71955338
$arr[3,1] = 98520
$arr[3,2] = "Ovaliderad"
$arr[3,3] = "LC"
$arr[3,4] = 222498
$arr[3,5] = "Blåsippa"
$arr[3,6] = "Hepatica nobilis"
$arr[3,7] = "Schreb."
$arr[3,8] = ""
$arr[3,9] = ""
$arr[3,10] = ""
$arr[3,11] = ""
$arr[3,12] = ""
$arr[3,13] = ""
$arr[3,14] = ""
$arr[3,15] = "Björkenäs, Gläntan, Upl"
$arr[3,16] = 688350.0772071786
$arr[3,17] = 6628595.795327608
$arr[3,18] = 5
$arr[3,19] = "Stockholm"
$arr[3,20] = "Norrtälje"
$arr[3,21] = "Uppland"
$arr[3,22] = "Rimbo"
$arr[3,23] = ""
$arr[3,24] = "2018-06-20"
$arr[3,25] = "00:00"
$arr[3,26] = "2018-06-20"
$arr[3,27] = "00:00"
$arr[3,28] = ""
$arr[3,29] = $false
$arr[3,30] = $false
$arr[3,31] = ""
$arr[3,32] = $false
$arr[3,33] = "Skogsmark"
$arr[3,34] = "Blandskog"
$arr[3,35] = ""
$arr[3,36] = ""
$arr[3,37] = ""
$arr[3,38] = ""
$arr[3,39] = ""
$arr[3,40] = ""
$arr[3,41] = ""
$arr[3,42] = ""
$arr[3,43] = ""
$arr[3,44] = ""
$arr[3,45] = ""
$arr[3,46] = ""
$arr[3,47] = ""
$arr[3,48] = "Beata Larsson"
$arr[3,49] = "Beata Larsson"
$arr[3,50] = ""

$arr[4,0] = 71955315
$arr[4,1] = 98520
$arr[4,2] = "Ovaliderad"
$arr[4,3] = "LC"
$arr[4,4] = 222498
$arr[4,5] = "Blåsippa"
$arr[4,6] = "Hepatica nobilis"
$arr[4,7] = "Schreb."
$arr[4,8] = ""
$arr[4,9] = ""
$arr[4,10] = ""
$arr[4,11] = ""
$arr[4,12] = ""
$arr[4,13] = ""
$arr[4,14] = ""
$arr[4,15] = "Björkenäs, Gläntan, Upl"
$arr[4,16] = 688375.1939352592
$arr[4,17] = 6628599.085711329
$arr[4,18] = 10
$arr[4,19] = "Stockholm"
$arr[4,20] = "Norrtälje"
$arr[4,21] = "Uppland"
$arr[4,22] = "Rimbo"
$arr[4,23] = ""
$arr[4,24] = "2018-06-20"
$arr[4,25] = "00:00"
$arr[4,26] = "2018-06-20"
$arr[4,27] = "00:00"
$arr[4,28] = ""
$arr[4,29] = $false
$arr[4,30] = $false
$arr[4,31] = ""
$arr[4,32] = $false
$arr[4,33] = "Skogsmark"
$arr[4,34] = "Blandskog"
$arr[4,35] = ""
$arr[4,36] = ""
$arr[4,37] = ""
$arr[4,38] = ""
$arr[4,39] = ""
$arr[4,40] = ""
$arr[4,41] = ""
$arr[4,42] = ""
$arr[4,43] = ""
$arr[4,44] = ""
$arr[4,45] = ""
$arr[4,46] = ""
$arr[4,47] = ""
$arr[4,48] = "Beata Larsson"
$arr[4,49] = "Beata Larsson"
$arr[4,50] = ""

$arr[5,0] = 71955359
$arr[5,1] = 98520
$arr[5,2] = "Ovaliderad"
$arr[5,3] = "LC"
$arr[5,4] = 222498
$arr[5,5] = "Blåsippa"
$arr[5,6] = "Hepatica nobilis"
$arr[5,7] = "Schreb."
$arr[5,8] = ""
$arr[5,9] = ""
$arr[5,10] = ""
$arr[5,11] = ""
$arr[5,12] = ""
$arr[5,13] = ""
$arr[5,14] = ""
$arr[5,15] = "Björkenäs, Gläntan, Upl"
$arr[5,16] = 688368.9025327084
$arr[5,17] = 6628573.527109651
$arr[5,18] = 5
$arr[5,19] = "Stockholm"
$arr[5,20] = "Norrtälje"
$arr[5,21] = "Uppland"
$arr[5,22] = "Rimbo"
$arr[5,23] = ""
$arr[5,24] = "2018-06-20"
$arr[5,25] = "00:00"
$arr[5,26] = "2018-06-20"
$arr[5,27] = "00:00"
$arr[5,28] = ""
$arr[5,29] = $false
$arr[5,30] = $false
$arr[5,31] = ""
$arr[5,32] = $false
$arr[5,33] = "Skogsmark"
$arr[5,34] = "Blandskog"
$arr[5,35] = ""
$arr[5,36] = ""
$arr[5,37] = ""
$arr[5,38] = ""
$arr[5,39] = ""
$arr[5,40] = ""
$arr[5,41] = ""
$arr[5,42] = ""
$arr[5,43] = ""
$arr[5,44] = ""
$arr[5,45] = ""
$arr[5,46] = ""
$arr[5,47] = ""
$arr[5,48] = "Beata Larsson"
$arr[5,49] = "Beata Larsson"
$arr[5,50] = ""

$arr[6,0] = 71955323
$arr[6,1] = 98520
$arr[6,2] = "Ovaliderad"
$arr[6,3] = "LC"
$arr[6,4] = 222498
$arr[6,5] = "Blåsippa"
$arr[6,6] = "Hepatica nobilis"
$arr[6,7] = "Schreb."
$arr[6,8] = ""
$arr[6,9] = ""
$arr[6,10] = ""
$arr[6,11] = ""
$arr[6,12] = ""
$arr[6,13] = ""
$arr[6,14] = ""
$arr[6,15] = "Björkenäs, Gläntan, Upl"
$arr[6,16] = 688371.3681368505
$arr[6,17] = 6628584.75759712
$arr[6,18] = 5
$arr[6,19] = "Stockholm"
$arr[6,20] = "Norrtälje"
$arr[6,21] = "Uppland"
$arr[6,22] = "Rimbo"
$arr[6,23] = ""
$arr[6,24] = "2018-06-20"
$arr[6,25] = "00:00"
$arr[6,26] = "2018-06-20"
$arr[6,27] = "00:00"
$arr[6,28] = ""
$arr[6,29] = $false
$arr[6,30] = $false
$arr[6,31] = ""
$arr[6,32] = $false
$arr[6,33] = "Skogsmark"
$arr[6,34] = "Blandskog"
$arr[6,35] = ""
$arr[6,36] = ""
$arr[6,37] = ""
$arr[6,38] = ""
$arr[6,39] = ""
$arr[6,40] = ""
$arr[6,41] = ""
$arr[6,42] = ""
$arr[6,43] = ""
$arr[6,44] = ""
$arr[6,45] = ""
$arr[6,46] = ""
$arr[6,47] = ""
$arr[6,48] = "Beata Larsson"
$arr[6,49] = "Beata Larsson"
$arr[6,50] = ""

$arr[7,0] = 71955445
$arr[7,1] = 98520
$arr[7,2] = "Ovaliderad"
$arr[7,3] = "LC"
$arr[7,4] = 222498
$arr[7,5] = "Blåsippa"
$arr[7,6] = "Hepatica nobilis"
$arr[7,7] = "Schreb."
$arr[7,8] = ""
$arr[7,9] = ""
$arr[7,10] = ""
$arr[7,11] = ""
$arr[7,12] = ""
$arr[7,13] = ""
$arr[7,14] = ""
$arr[7,15] = "Björkenäs, Gläntan, Upl"
$arr[7,16] = 688339.897607106
$arr[7,17] = 6628627.082556061
$arr[7,18] = 5
$arr[7,19] = "Stockholm"
$arr[7,20] = "Norrtälje"
$arr[7,21] = "Uppland"
$arr[7,22] = "Rimbo"
$arr[7,23] = ""
$arr[7,24] = "2018-06-20"
$arr[7,25] = "00:00"
$arr[7,26] = "2018-06-20"
$arr[7,27] = "00:00"
$arr[7,28] = ""
$arr[7,29] = $false
$arr[7,30] = $false
$arr[7,31] = ""
$arr[7,32] = $false
$arr[7,33] = "Skogsmark"
$arr[7,34] = "Blandskog"
$arr[7,35] = ""
$arr[7,36] = ""
$arr[7,37] = ""
$arr[7,38] = ""
$arr[7,39] = ""
$arr[7,40] = ""
$arr[7,41] = ""
$arr[7,42] = ""
$arr[7,43] = ""
$arr[7,44] = ""
$arr[7,45] = ""
$arr[7,46] = ""
$arr[7,47] = ""
$arr[7,48] = "Beata Larsson"
$arr[7,49] = "Beata Larsson"
$arr[7,50] = ""

$arr[8,0] = 71965499
$arr[8,1] = 101120
$arr[8,2] = "Ovaliderad"
$arr[8,3] = "LC"
$arr[8,4] = 222002
$arr[8,5] = "Underviol"
$arr[8,6] = "Viola mirabilis"
$arr[8,7] = "L."
$arr[8,8] = ""
$arr[8,9] = ""
$arr[8,10] = ""
$arr[8,11] = ""
$arr[8,12] = ""
$arr[8,13] = ""
$arr[8,14] = ""
$arr[8,15] = "Björkenäs, Gläntan, Upl"
$arr[8,16] = 688317.766267107
$arr[8,17] = 6628654.735997377
$arr[8,18] = 5
$arr[8,19] = "Stockholm"
$arr[8,20] = "Norrtälje"
$arr[8,21] = "Uppland"
$arr[8,22] = "Rimbo"
$arr[8,23] = ""
$arr[8,24] = "2018-06-21"
$arr[8,25] = "00:00"
$arr[8,26] = "2018-06-21"
$arr[8,27] = "00:00"
$arr[8,28] = ""
$arr[8,29] = $false
$arr[8,30] = $false
$arr[8,31] = ""
$arr[8,32] = $false
$arr[8,33] = "Skogsmark"
$arr[8,34] = "Blandskog"
$arr[8,35] = ""
$arr[8,36] = ""
$arr[8,37] = ""
$arr[8,38] = ""
$arr[8,39] = ""
$arr[8,40] = ""
$arr[8,41] = ""
$arr[8,42] = ""
$arr[8,43] = ""
$arr[8,44] = ""
$arr[8,45] = ""
$arr[8,46] = ""
$arr[8,47] = ""
$arr[8,48] = "Beata Larsson"
$arr[8,49] = "Beata Larsson"
$arr[8,50] = ""

$arr[9,0] = 71964226
$arr[9,1] = 98520
$arr[9,2] = "Ovaliderad"
$arr[9,3] = "LC"
$arr[9,4] = 222498
$arr[9,5] = "Blåsippa"
$arr[9,6] = "Hepatica nobilis"
$arr[9,7] = "Schreb."
$arr[9,8] = ""
$arr[9,9] = ""
$arr[9,10] = ""
$arr[9,11] = ""
$arr[9,12] = ""
$arr[9,13] = ""
$arr[9,14] = ""
$arr[9,15] = "Björkenäs, Gläntan, Upl"
$arr[9,16] = 688316.8952691283
$arr[9,17] = 6628631.976152469
$arr[9,18] = 5
$arr[9,19] = "Stockholm"
$arr[9,20] = "Norrtälje"
$arr[9,21] = "Uppland"
$arr[9,22] = "Rimbo"
$arr[9,23] = ""
$arr[9,24] = "2018-06-21"
$arr[9,25] = "00:00"
$arr[9,26] = "2018-06-21"
$arr[9,27] = "00:00"
$arr[9,28] = ""
$arr[9,29] = $false
$arr[9,30] = $false
$arr[9,31] = ""
$arr[9,32] = $false
$arr[9,33] = "Skogsmark"
$arr[9,34] = "Blandskog"
$arr[9,35] = ""
$arr[9,36] = ""
$arr[9,37] = ""
$arr[9,38] = ""
$arr[9,39] = ""
$arr[9,40] = ""
$arr[9,41] = ""
$arr[9,42] = ""
$arr[9,43] = ""
$arr[9,44] = ""
$arr[9,45] = ""
$arr[9,46] = ""
$arr[9,47] = ""
$arr[9,48] = "Beata Larsson"
$arr[9,49] = "Beata Larsson"
$arr[9,50] = ""

$arr[10,0] = 71965450
$arr[10,1] = 98520
$arr[10,2] = "Ovaliderad"
$arr[10,3] = "LC"
$arr[10,4] = 222498
$arr[10,5] = "Blåsippa"
$arr[10,6] = "Hepatica nobilis"
$arr[10,7] = "Schreb."
$arr[10,8] = ""
$arr[10,9] = ""
$arr[10,10] = ""
$arr[10,11] = ""
$arr[10,12] = ""
$arr[10,13] = ""
$arr[10,14] = ""
$arr[10,15] = "Björkenäs, Gläntan, Upl"
$arr[10,16] = 688313.4256247933
$arr[10,17] = 6628660.573880369
$arr[10,18] = 5
$arr[10,19] = "Stockholm"
$arr[10,20] = "Norrtälje"
$arr[10,21] = "Uppland"
$arr[10,22] = "Rimbo"
$arr[10,23] = ""
$arr[10,24] = "2018-06-21"
$arr[10,25] = "00:00"
$arr[10,26] = "2018-06-21"
$arr[10,27] = "00:00"
$arr[10,28] = ""
$arr[10,29] = $false
$arr[10,30] = $false
$arr[10,31] = ""
$arr[10,32] = $false
$arr[10,33] = "Skogsmark"
$arr[10,34] = "Blandskog"
$arr[10,35] = ""
$arr[10,36] = ""
$arr[10,37] = ""
$arr[10,38] = ""
$arr[10,39] = ""
$arr[10,40] = ""
$arr[10,41] = ""
$arr[10,42] = ""
$arr[10,43] = ""
$arr[10,44] = ""
$arr[10,45] = ""
$arr[10,46] = ""
$arr[10,47] = ""
$arr[10,48] = "Beata Larsson"
$arr[10,49] = "Beata Larsson"
$arr[10,50] = ""

$arr[11,0] = 71964234
$arr[11,1] = 101120
$arr[11,2] = "Ovaliderad"
$arr[11,3] = "LC"
$arr[11,4] = 222002
$arr[11,5] = "Underviol"
$arr[11,6] = "Viola mirabilis"
$arr[11,7] = "L."
$arr[11,8] = ""
$arr[11,9] = ""
$arr[11,10] = ""
$arr[11,11] = ""
$arr[11,12] = ""
$arr[11,13] = ""
$arr[11,14] = ""
$arr[11,15] = "Björkenäs, Gläntan, Upl"
$arr[11,16] = 688313.8690585461
$arr[11,17] = 6628631.823030195
$arr[11,18] = 5
$arr[11,19] = "Stockholm"
$arr[11,20] = "Norrtälje"
$arr[11,21] = "Uppland"
$arr[11,22] = "Rimbo"
$arr[11,23] = ""
$arr[11,24] = "2018-06-21"
$arr[11,25] = "00:00"
$arr[11,26] = "2018-06-21"
$arr[11,27] = "00:00"
$arr[11,28] = ""
$arr[11,29] = $false
$arr[11,30] = $false
$arr[11,31] = ""
$arr[11,32] = $false
$arr[11,33] = "Skogsmark"
$arr[11,34] = "Blandskog"
$arr[11,35] = ""
$arr[11,36] = ""
$arr[11,37] = ""
$arr[11,38] = ""
$arr[11,39] = ""
$arr[11,40] = ""
$arr[11,41] = ""
$arr[11,42] = ""
$arr[11,43] = ""
$arr[11,44] = ""
$arr[11,45] = ""
$arr[11,46] = ""
$arr[11,47] = ""
$arr[11,48] = "Beata Larsson"
$arr[11,49] = "Beata Larsson"
$arr[11,50] = ""

$arr[12,0] = 71965509
$arr[12,1] = 98520
$arr[12,2] = "Ovaliderad"
$arr[12,3] = "LC"
$arr[12,4] = 222498
$arr[12,5] = "Blåsippa"
$arr[12,6] = "Hepatica nobilis"
$arr[12,7] = "Schreb."
$arr[12,8] = ""
$arr[12,9] = ""
$arr[12,10] = ""
$arr[12,11] = ""
$arr[12,12] = ""
$arr[12,13] = ""
$arr[12,14] = ""
$arr[12,15] = "Björkenäs, Gläntan, Upl"
$arr[12,16] = 688317.766267107
$arr[12,17] = 6628654.735997377
$arr[12,18] = 5
$arr[12,19] = "Stockholm"
$arr[12,20] = "Norrtälje"
$arr[12,21] = "Uppland"
$arr[12,22] = "Rimbo"
$arr[12,23] = ""
$arr[12,24] = "2018-06-21"
$arr[12,25] = "00:00"
$arr[12,26] = "2018-06-21"
$arr[12,27] = "00:00"
$arr[12,28] = ""
$arr[12,29] = $false
$arr[12,30] = $false
$arr[12,31] = ""
$arr[12,32] = $false
$arr[12,33] = "Skogsmark"
$arr[12,34] = "Blandskog"
$arr[12,35] = ""
$arr[12,36] = ""
$arr[12,37] = ""
$arr[12,38] = ""
$arr[12,39] = ""
$arr[12,40] = ""
$arr[12,41] = ""
$arr[12,42] = ""
$arr[12,43] = ""
$arr[12,44] = ""
$arr[12,45] = ""
$arr[12,46] = ""
$arr[12,47] = ""
$arr[12,48] = "Beata Larsson"
$arr[12,49] = "Beata Larsson"
$arr[12,50] = ""

$arr[13,0] = 71964207
$arr[13,1] = 98431
$arr[13,2] = "Ovaliderad"
$arr[13,3] = "LC"
$arr[13,4] = 222771
$arr[13,5] = "Svart trolldruva"
$arr[13,6] = "Actaea spicata"
$arr[13,7] = "L."
$arr[13,8] = ""
$arr[13,9] = ""
$arr[13,10] = ""
$arr[13,11] = ""
$arr[13,12] = ""
$arr[13,13] = ""
$arr[13,14] = ""
$arr[13,15] = "Björkenäs, Gläntan, Upl"
$arr[13,16] = 688318.7650556953
$arr[13,17] = 6628625.003604913
$arr[13,18] = 5
$arr[13,19] = "Stockholm"
$arr[13,20] = "Norrtälje"
$arr[13,21] = "Uppland"
$arr[13,22] = "Rimbo"
$arr[13,23] = ""
$arr[13,24] = "2018-06-21"
$arr[13,25] = "00:00"
$arr[13,26] = "2018-06-21"
$arr[13,27] = "00:00"
$arr[13,28] = ""
$arr[13,29] = $false
$arr[13,30] = $false
$arr[13,31] = ""
$arr[13,32] = $false
$arr[13,33] = "Skogsmark"
$arr[13,34] = "Blandskog"
$arr[13,35] = ""
$arr[13,36] = ""
$arr[13,37] = ""
$arr[13,38] = ""
$arr[13,39] = ""
$arr[13,40] = ""
$arr[13,41] = ""
$arr[13,42] = ""
$arr[13,43] = ""
$arr[13,44] = ""
$arr[13,45] = ""
$arr[13,46] = ""
$arr[13,47] = ""
$arr[13,48] = "Beata Larsson"
$arr[13,49] = "Beata Larsson"
$arr[13,50] = ""

$arr[14,0] = 71966065
$arr[14,1] = 98520
$arr[14,2] = "Ovaliderad"
$arr[14,3] = "LC"
$arr[14,4] = 222498
$arr[14,5] = "Blåsippa"
$arr[14,6] = "Hepatica nobilis"
$arr[14,7] = "Schreb."
$arr[14,8] = ""
$arr[14,9] = ""
$arr[14,10] = ""
$arr[14,11] = ""
$arr[14,12] = ""
$arr[14,13] = ""
$arr[14,14] = ""
$arr[14,15] = "Björkenäs, Gläntan, Upl"
$arr[14,16] = 688353.4294477219
$arr[14,17] = 6628599.498588714
$arr[14,18] = 5
$arr[14,19] = "Stockholm"
$arr[14,20] = "Norrtälje"
$arr[14,21] = "Uppland"
$arr[14,22] = "Rimbo"
$arr[14,23] = ""
$arr[14,24] = "2018-06-21"
$arr[14,25] = "00:00"
$arr[14,26] = "2018-06-21"
$arr[14,27] = "00:00"
$arr[14,28] = ""
$arr[14,29] = $false
$arr[14,30] = $false
$arr[14,31] = ""
$arr[14,32] = $false
$arr[14,33] = "Skogsmark"
$arr[14,34] = "Blandskog"
$arr[14,35] = ""
$arr[14,36] = ""
$arr[14,37] = ""
$arr[14,38] = ""
$arr[14,39] = ""
$arr[14,40] = ""
$arr[14,41] = ""
$arr[14,42] = ""
$arr[14,43] = ""
$arr[14,44] = ""
$arr[14,45] = ""
$arr[14,46] = ""
$arr[14,47] = ""
$arr[14,48] = "Beata Larsson"
$arr[14,49] = "Beata Larsson"
$arr[14,50] = ""

$arr[15,0] = 71965472
$arr[15,1] = 101120
$arr[15,2] = "Ovaliderad"
$arr[15,3] = "LC"
$arr[15,4] = 222002
$arr[15,5] = "Underviol"
$arr[15,6] = "Viola mirabilis"
$arr[15,7] = "L."
$arr[15,8] = ""
$arr[15,9] = ""
$arr[15,10] = ""
$arr[15,11] = ""
$arr[15,12] = ""
$arr[15,13] = ""
$arr[15,14] = ""
$arr[15,15] = "Björkenäs, Gläntan, Upl"
$arr[15,16] = 688325.5355727162
$arr[15,17] = 6628651.090780158
$arr[15,18] = 5
$arr[15,19] = "Stockholm"
$arr[15,20] = "Norrtälje"
$arr[15,21] = "Uppland"
$arr[15,22] = "Rimbo"
$arr[15,23] = ""
$arr[15,24] = "2018-06-21"
$arr[15,25] = "00:00"
$arr[15,26] = "2018-06-21"
$arr[15,27] = "00:00"
$arr[15,28] = ""
$arr[15,29] = $false
$arr[15,30] = $false
$arr[15,31] = ""
$arr[15,32] = $false
$arr[15,33] = "Skogsmark"
$arr[15,34] = "Blandskog"
$arr[15,35] = ""
$arr[15,36] = ""
$arr[15,37] = ""
$arr[15,38] = ""
$arr[15,39] = ""
$arr[15,40] = ""
$arr[15,41] = ""
$arr[15,42] = ""
$arr[15,43] = ""
$arr[15,44] = ""
$arr[15,45] = ""
$arr[15,46] = ""
$arr[15,47] = ""
$arr[15,48] = "Beata Larsson"
$arr[15,49] = "Beata Larsson"
$arr[15,50] = ""

$arr[16,0] = 73101994
$arr[16,1] = 98520
$arr[16,2] = "Ovaliderad"
$arr[16,3] = "LC"
$arr[16,4] = 222498
$arr[16,5] = "Blåsippa"
$arr[16,6] = "Hepatica nobilis"
$arr[16,7] = "Schreb."
$arr[16,8] = ""
$arr[16,9] = ""
$arr[16,10] = ""
$arr[16,11] = ""
$arr[16,12] = ""
$arr[16,13] = ""
$arr[16,14] = ""
$arr[16,15] = "Björkenäs, Upl"
$arr[16,16] = 688330.0904902341
$arr[16,17] = 6628621.033502498
$arr[16,18] = 10
$arr[16,19] = "Stockholm"
$arr[16,20] = "Norrtälje"
$arr[16,21] = "Uppland"
$arr[16,22] = "Rimbo"
$arr[16,23] = ""
$arr[16,24] = "2018-06-26"
$arr[16,25] = "00:00"
$arr[16,26] = "2018-06-26"
$arr[16,27] = "00:00"
$arr[16,28] = ""
$arr[16,29] = $false
$arr[16,30] = $false
$arr[16,31] = ""
$arr[16,32] = $false
$arr[16,33] = ""
$arr[16,34] = ""
$arr[16,35] = ""
$arr[16,36] = ""
$arr[16,37] = ""
$arr[16,38] = ""
$arr[16,39] = ""
$arr[16,40] = ""
$arr[16,41] = ""
$arr[16,42] = ""
$arr[16,43] = ""
$arr[16,44] = ""
$arr[16,45] = ""
$arr[16,46] = ""
$arr[16,47] = ""
$arr[16,48] = "Niina Sallmén"
$arr[16,49] = "Niina Sallmén"
$arr[16,50] = ""

$arr[17,0] = 92897138
$arr[17,1] = 56521
$arr[17,2] = "Ovaliderad"
$arr[17,3] = "NT"
$arr[17,4] = 103035
$arr[17,5] = "Kråka"
$arr[17,6] = "Corvus corone"
$arr[17,7] = "Linnaeus, 1758"
$arr[17,8] = ""
$arr[17,9] = ""
$arr[17,10] = ""
$arr[17,11] = ""
$arr[17,12] = ""
$arr[17,13] = ""
$arr[17,14] = ""
$arr[17,15] = "Bålbroskogen, Upl"
$arr[17,16] = 688356.6664694136
$arr[17,17] = 6628485.575551489
$arr[17,18] = 50
$arr[17,19] = "Stockholm"
$arr[17,20] = "Norrtälje"
$arr[17,21] = "Uppland"
$arr[17,22] = "Rimbo"
$arr[17,23] = ""
$arr[17,24] = "2020-05-14"
$arr[17,25] = "00:00"
$arr[17,26] = "2020-05-14"
$arr[17,27] = "00:00"
$arr[17,28] = ""
$arr[17,29] = $false
$arr[17,30] = $false
$arr[17,31] = ""
$arr[17,32] = $false
$arr[17,33] = ""
$arr[17,34] = ""
$arr[17,35] = ""
$arr[17,36] = ""
$arr[17,37] = ""
$arr[17,38] = ""
$arr[17,39] = ""
$arr[17,40] = ""
$arr[17,41] = ""
$arr[17,42] = ""
$arr[17,43] = ""
$arr[17,44] = ""
$arr[17,45] = ""
$arr[17,46] = ""
$arr[17,47] = ""
$arr[17,48] = "Mattias Lif"
$arr[17,49] = "Mattias Lif"
$arr[17,50] = ""

$arr[18,0] = 88441956
$arr[18,1] = 85318
$arr[18,2] = "Ovaliderad"
$arr[18,3] = "VU"
$arr[18,4] = 3767
$arr[18,5] = "Violettfläckig spindling"
$arr[18,6] = "Cortinarius violaceomaculatus"
$arr[18,7] = "Brandrud"
$arr[18,8] = "6"
$arr[18,9] = "fruktkroppar"
$arr[18,10] = ""
$arr[18,11] = ""
$arr[18,12] = ""
$arr[18,13] = ""
$arr[18,14] = ""
$arr[18,15] = "Rimbo N, Upl"
$arr[18,16] = 688413.8515789268
$arr[18,17] = 6628464.743477113
$arr[18,18] = 25
$arr[18,19] = "Stockholm"
$arr[18,20] = "Norrtälje"
$arr[18,21] = "Uppland"
$arr[18,22] = "Rimbo"
$arr[18,23] = ""
$arr[18,24] = "2020-10-06"
$arr[18,25] = "00:00"
$arr[18,26] = "2020-10-06"
$arr[18,27] = "00:00"
$arr[18,28] = ""
$arr[18,29] = $false
$arr[18,30] = $false
$arr[18,31] = ""
$arr[18,32] = $false
$arr[18,33] = ""
$arr[18,34] = ""
$arr[18,35] = ""
$arr[18,36] = ""
$arr[18,37] = ""
$arr[18,38] = ""
$arr[18,39] = ""
$arr[18,40] = ""
$arr[18,41] = ""
$arr[18,42] = ""
$arr[18,43] = ""
$arr[18,44] = ""
$arr[18,45] = ""
$arr[18,46] = ""
$arr[18,47] = ""
$arr[18,48] = "Ossian Rydebjörk"
$arr[18,49] = "Ossian Rydebjörk, Birgitta Wasstorp"
$arr[18,50] = "Svampar i Roslagen"

$rng = $ws.Range("A2:AY20")
$rng.Value2 = $arr

Write-Host "Done applying row rotation."